$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mato Grosso"
$ws.Range("B2").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C2").Value = 2.691531905194921

$ws.Range("A3").Value = "Goiás"
$ws.Range("B3").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C3").Value = 1.501862535602321

$ws.Range("B4").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C4").Value = 1.443172347238672

$ws.Range("A5").Value = "Paraíba"
$ws.Range("B5").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C5").Value = 0.8595870483680415

$ws.Range("A6").Value = "Minas Gerais"
$ws.Range("B6").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C6").Value = 0.5646447141197939

$ws.Range("A7").Value = "São Paulo"
$ws.Range("B7").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C7").Value = 0.5151860777804842

$ws.Range("B8").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C8").Value = -1.418749418647174
$ws.Range("D8").Value = "27º"

$ws.Range("B9").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C9").Value = -0.2144178582346825

$ws.Range("B10").Value = "Diferença 2023/04 - 2022/04"
$ws.Range("C10").Value = 0.1843642635987592
